$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Session 9 (col O) and Session 10 (col P) attendance marks for rows 7..86
$oVals = @("A", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P")
$pVals = @("P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "A", "P", "P", "P", "P", "P", "A", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "A", "P", "P", "P", "P", "P", "P", "A", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "P", "A", "P", "P", "P", "P", "P", "A", "P", "P", "P", "P")

for ($i = 0; $i -lt $oVals.Length; $i++) {
    $r = 7 + $i
    # New cells inherit the bordered/centered style used by the other attendance columns (copy from column H)
    $ws.Range("H" + $r).Copy() | Out-Null
    $ws.Range("O" + $r).PasteSpecial(-4122) | Out-Null
    $ws.Range("H" + $r).Copy() | Out-Null
    $ws.Range("P" + $r).PasteSpecial(-4122) | Out-Null

    $ws.Range("O" + $r).Value = $oVals[$i]
    $ws.Range("P" + $r).Value = $pVals[$i]
}
